$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-5 (Ifne -> Ifnar1 pairs)

# Row 2 (MuSCs -> ECs)
$ws.Range("G2").Value = 0.057308
$ws.Range("H2").Value = 0.171924
$ws.Range("M2").Value = 21.602164
$ws.Range("N2").Value = 64.80649199999999
$ws.Range("O2").Value = 0.2504461301095448
$ws.Range("P2").Value = 0.2504461301095448
$ws.Range("Q2").Value = 1.237976814512
$ws.Range("R2").Value = 11.141791330608
$ws.Range("S2").Value = 0.2504461301095448
$ws.Range("T2").Value = 0.2504461301095448

# Row 3 (MuSCs -> FAPs)
$ws.Range("G3").Value = 0.057308
$ws.Range("H3").Value = 0.171924
$ws.Range("O3").Value = 0.2379112932771326
$ws.Range("P3").Value = 0.2379112932771326
$ws.Range("Q3").Value = 1.1760160353
$ws.Range("R3").Value = 10.5841443177
$ws.Range("S3").Value = 0.2379112932771326
$ws.Range("T3").Value = 0.2379112932771326

# Row 4 (MuSCs -> MuSCs)
$ws.Range("G4").Value = 0.057308
$ws.Range("H4").Value = 0.171924
$ws.Range("M4").Value = 11.57857933333334
$ws.Range("N4").Value = 34.735738
$ws.Range("O4").Value = 0.1342370322806403
$ws.Range("P4").Value = 0.1342370322806404
$ws.Range("Q4").Value = 0.6635452244346668
$ws.Range("R4").Value = 5.971907019912001
$ws.Range("S4").Value = 0.1342370322806403
$ws.Range("T4").Value = 0.1342370322806404

# Row 5 (MuSCs -> Resolving-Mac)
$ws.Range("G5").Value = 0.057308
$ws.Range("H5").Value = 0.171924
$ws.Range("M5").Value = 32.55301433333333
$ws.Range("N5").Value = 97.659043
$ws.Range("O5").Value = 0.3774055443326823
$ws.Range("P5").Value = 0.3774055443326824
$ws.Range("Q5").Value = 1.865548145414666
$ws.Range("R5").Value = 16.789933308732
$ws.Range("S5").Value = 0.3774055443326823
$ws.Range("T5").Value = 0.3774055443326824
